# This script rewrites the comp_loinc workbook structure:
#  - Renames all 12 sheets (tab order stays 1..12, only names + contents change)
#  - Replaces the header row (row 1) of each sheet with the new set of columns
#
# Final tab order / names:
#  1  Thing
#  2  Loinc
#  3  PartClass
#  4  ComponentClass
#  5  SystemClass
#  6  MethodClass
#  7  TimeClass
#  8  PropertyClass
#  9  ScaleClass
#  10 LoincCodeClass
#  11 LoincCodeClassIntersection
#  12 LoincCodeClassNonIntersection

$wb = $excel.ActiveWorkbook

# Column header sets used by the various sheets after the rewrite.

$thingCols = @("id","label","description","subClassOf","equivalentClasses")

$loincCols = @("codes","parts")

$partLikeCols = @("subClassOf","part_number","part_type","part_name","part_display_name","id","label","description","equivalentClasses")

$loincCodeCols = @(
    "loinc_number","long_common_name","formal_name","short_name","status","loinc_class","loinc_class_type",
    "has_component","has_property","has_time","has_system","has_scale","has_method",
    "has_component_analyte","has_component_challenge","has_component_count","has_component_adjustment",
    "has_time_core","has_time_modifier","has_system_core","has_system_super_system",
    "semantic_analyte_gene","syntax_analyte_core","syntax_analyte_suffix","syntax_analyte_divisor",
    "syntax_analyte_divisor_suffix","syntax_analyte_numerator","id","label","description","subClassOf","equivalentClasses"
)

# Ordered list describing each tab's final name and which header set it gets.
$sheetPlan = @(
    @{ Name = "Thing";                          Cols = $thingCols },
    @{ Name = "Loinc";                           Cols = $loincCols },
    @{ Name = "PartClass";                       Cols = $partLikeCols },
    @{ Name = "ComponentClass";                  Cols = $partLikeCols },
    @{ Name = "SystemClass";                     Cols = $partLikeCols },
    @{ Name = "MethodClass";                     Cols = $partLikeCols },
    @{ Name = "TimeClass";                       Cols = $partLikeCols },
    @{ Name = "PropertyClass";                   Cols = $partLikeCols },
    @{ Name = "ScaleClass";                      Cols = $partLikeCols },
    @{ Name = "LoincCodeClass";                  Cols = $loincCodeCols },
    @{ Name = "LoincCodeClassIntersection";      Cols = $loincCodeCols },
    @{ Name = "LoincCodeClassNonIntersection";   Cols = $loincCodeCols }
)

# Renaming straight to the final names can collide with another sheet that
# still carries its old name (e.g. tab 3 needs to become "PartClass" while
# tab 5 is still named "PartClass" until it's processed). Avoid that by
# first moving every sheet to a unique scratch name, then assigning the
# real final names in a second pass.
for ($i = 0; $i -lt $sheetPlan.Length; $i++) {
    $ws = $wb.Worksheets.Item($i + 1)
    $ws.Name = "__scratch_$i"
}

for ($i = 0; $i -lt $sheetPlan.Length; $i++) {
    $plan = $sheetPlan[$i]
    $ws = $wb.Worksheets.Item($i + 1)

    # Clear any existing content first so stale cells beyond the new header
    # row/columns don't linger.
    $ws.Cells.Clear()

    $ws.Name = $plan.Name

    $cols = $plan.Cols
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $ws.Cells.Item(1, $c + 1).Value = $cols[$c]
    }
}

Write-Host "rewrote $($sheetPlan.Length) sheets"
